# Fix Training Data Issue (#48)
# The "Date" column (BF) held the source filename-derived text
# "5-30-2013-14" for every row. Re-express it as the real game date
# "2014-05-30" (the NBA stats page date was off by a day relative to
# how the file was named).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
# Force plain text storage so Excel doesn't reinterpret the
# "YYYY-MM-DD"-shaped string as a date serial number.
$dateRange.NumberFormat = "@"
$dateRange.Value = "2014-05-30"
